# Append four more entries (rows 8-11) to Sheet1, duplicating existing
# dogs (Moss, Arlo, Wern Joe, Kinloch Heath) as new entries, then update
# the active selection to reflect where the user ended up afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRows = @(
    @{ Row = 8;  A = "Moss";          B = "Annie  Vanderlinck"; C = "D. T. Edwards";   D = 39861; E = 123435;      F = "Mirk";           G = "Llangwn Lass";  H = "PicturesOrig/AnnieVanderlinck_Moss.jpg" },
    @{ Row = 9;  A = "Arlo";          B = "Folke Noertemann";   C = "Jean-Luc Censier"; D = 45433; E = 12345;       F = "Noi";            G = "Penny";         H = "PicturesOrig/FolkeNoertemann_Arlo.jpg" },
    @{ Row = 10; A = "Wern Joe";      B = "Folke Noertemann";   C = "M.D. Jones";       D = 44247; E = "1112355 HG"; F = "Kinloch Sweep";  G = "crazy mommy";   H = "PicturesOrig/FolkeNoertemann_Joe.jpg" },
    @{ Row = 11; A = "Kinloch Heath"; B = "Folke Noertemann";   C = "Angie Driscoll";   D = 43322; E = 33333;       F = "no idea";        G = "Kinloch Penny"; H = "PicturesOrig/FolkeNoertemann_Heath.jpg" }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value2 = $entry.A
    $ws.Cells.Item($r, 2).Value2 = $entry.B
    $ws.Cells.Item($r, 3).Value2 = $entry.C
    $ws.Cells.Item($r, 4).Value2 = $entry.D
    $ws.Cells.Item($r, 5).Value2 = $entry.E
    $ws.Cells.Item($r, 6).Value2 = $entry.F
    $ws.Cells.Item($r, 7).Value2 = $entry.G
    $ws.Cells.Item($r, 8).Value2 = $entry.H
}

$ws.Range("B21").Select()
